$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'69.334.05"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +1.27%  '
$ws.Range('D3').Value = "'3.944.98"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +0.26%  '
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').Value = "'494.97"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.27%  '
$ws.Range('D6').Value = "'147.88"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.31%  '
$ws.Range('D7').Value = "'0.625"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.95%  '
$ws.Range('E8').Value = '  +0.05%  '
$ws.Range('D9').Value = "'0.733"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.42%  '
$ws.Range('E10').Value = '  +4.42%  '
$ws.Range('E11').Value = '  -0.73%  '
$ws.Range('D12').Value = "'43.42"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +0.73%  '
$ws.Range('D13').Value = "'10.48"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -2.27%  '
$ws.Range('D14').Value = "'4.575.29"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +0.20%  '
$ws.Range('D15').Value = "'3.953.25"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +0.48%  '
$ws.Range('D16').Value = "'14.27"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -3.17%  '
$ws.Range('E17').Value = '  -0.32%  '
$ws.Range('E18').Value = '  +4.75%  '
$ws.Range('D19').Value = "'19.96"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.37%  '
$ws.Range('D20').Value = "'69.364.27"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +1.16%  '
$ws.Range('D21').Value = "'439.32"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -1.18%  '
$ws.Range('D22').Value = "'3.45"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -1.47%  '
$ws.Range('D23').Value = "'14.66"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -1.63%  '
$ws.Range('D24').Value = "'88.94"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.21%  '
$ws.Range('D25').Value = "'12.05"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +6.84%  '
$ws.Range('E26').Value = '  +4.01%  '
$ws.Range('E27').Value = '  -3.36%  '
$ws.Range('D28').Value = "'37.23"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -4.51%  '
$ws.Range('D29').Value = "'5.67"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -3.12%  '
$ws.Range('D30').Value = "'701.67"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -2.47%  '
$ws.Range('D31').Value = "'13.40"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -2.27%  '
$ws.Range('E32').Value = '  -0.27%  '
$ws.Range('E33').Value = '  -0.67%  '
$ws.Range('D34').Value = "'0.463"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +15.15%  '
$ws.Range('D35').Value = "'0.0₃0894"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -2.55%  '
$ws.Range('D36').Value = "'62.61"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +2.79%  '
$ws.Range('E37').Value = '  -2.47%  '
$ws.Range('D38').Value = "'41.11"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -2.50%  '
$ws.Range('E39').Value = '  +0.71%  '
$ws.Range('D40').Value = "'0.999"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -0.12%  '
$ws.Range('E41').Value = '  -0.02%  '
$ws.Range('D42').Value = "'0.0490"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +1.20%  '
$ws.Range('D43').Value = "'2.92"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -2.13%  '
$ws.Range('D44').Value = "'3.09"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -1.93%  '
$ws.Range('E45').Value = '  +1.94%  '
$ws.Range('E46').Value = '  +0.82%  '
$ws.Range('E47').Value = '  +6.98%  '
$ws.Range('E48').Value = '  +5.57%  '
$ws.Range('D49').Value = "'3.40"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -1.03%  '
$ws.Range('B50').Value = 'BabyDogeCoin'
$ws.Range('C50').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D50').Value = "'0.0₆0346"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +1.00%  '
$ws.Range('B51').Value = 'ARBITRUM'
$ws.Range('C51').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D51').Value = "'2.10"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -2.77%  '
